# Apply crypto price/volume updates per the commit diff.
# (cells under column D whose new text looks like a plain number are
#  written with a leading apostrophe so Excel keeps them as literal text,
#  matching the original inlineStr cells instead of coercing to a Number
#  and silently dropping things like trailing zeros.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.804.51'
$ws.Range('E2').Value = '  +0.69%  '
$ws.Range('D3').Value = '1.703.07'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('D4').Value = '''1.003'
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').Value = '''316.49'
$ws.Range('E5').Value = '  -0.10%  '
$ws.Range('D6').Value = '''1.004'
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('D7').Value = '''0.3934'
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('D9').Value = '''1.518'
$ws.Range('E9').Value = '  -1.52%  '
$ws.Range('D10').Value = '''1.003'
$ws.Range('E10').Value = '  +0.24%  '
$ws.Range('D11').Value = '''53.45'
$ws.Range('E11').Value = '  -2.28%  '
$ws.Range('D12').Value = '''0.08900'
$ws.Range('E12').Value = '  +0.97%  '
$ws.Range('D13').Value = '''7.310'
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('D14').Value = '''23.70'
$ws.Range('E14').Value = '  +1.02%  '
$ws.Range('D15').Value = '''8.021'
$ws.Range('E15').Value = '  +4.80%  '
$ws.Range('D16').Value = '''0.00001328'
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('D17').Value = '1.714.39'
$ws.Range('E17').Value = '  +0.95%  '
$ws.Range('D18').Value = '''100.43'
$ws.Range('E18').Value = '  -0.88%  '
$ws.Range('D19').Value = '''0.07042'
$ws.Range('E19').Value = '  -0.78%  '
$ws.Range('D20').Value = '''19.74'
$ws.Range('E20').Value = '  -0.32%  '
$ws.Range('D21').Value = '''7.082'
$ws.Range('E21').Value = '  +2.52%  '
$ws.Range('D22').Value = '''1.002'
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('D23').Value = '''14.54'
$ws.Range('E23').Value = '  +2.69%  '
$ws.Range('D24').Value = '24.781.62'
$ws.Range('E24').Value = '  +0.66%  '
$ws.Range('D25').Value = '''3.236'
$ws.Range('E25').Value = '  +3.24%  '
$ws.Range('D26').Value = '''2.360'
$ws.Range('E26').Value = '  +1.07%  '
$ws.Range('D27').Value = '''22.83'
$ws.Range('E27').Value = '  +1.62%  '
$ws.Range('D28').Value = '''162.55'
$ws.Range('E28').Value = '  +1.61%  '
$ws.Range('D29').Value = '''8.483'
$ws.Range('E29').Value = '  +11.08%  '
$ws.Range('D30').Value = '''136.85'
$ws.Range('E30').Value = '  +1.91%  '
$ws.Range('D31').Value = '''5.179'
$ws.Range('E31').Value = '  -1.29%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '''7.665'
$ws.Range('E32').Value = '  +2.87%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '''0.08900'
$ws.Range('E33').Value = '  +3.94%  '
$ws.Range('D34').Value = '''1.085'
$ws.Range('E34').Value = '  -2.39%  '
$ws.Range('E35').Value = '  -1.53%  '
$ws.Range('D36').Value = '''1.983'
$ws.Range('E36').Value = '  +1.95%  '
$ws.Range('D37').Value = '''0.2761'
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('E38').Value = '  -1.68%  '
$ws.Range('D39').Value = '''0.09211'
$ws.Range('E39').Value = '  +0.97%  '
$ws.Range('D40').Value = '''0.02779'
$ws.Range('E40').Value = '  -0.76%  '
$ws.Range('D41').Value = '''1.464'
$ws.Range('E41').Value = '  -0.21%  '
$ws.Range('D42').Value = '''0.7714'
$ws.Range('E42').Value = '  -0.74%  '
$ws.Range('E43').Value = '  +2.97%  '
$ws.Range('D44').Value = '''0.7200'
$ws.Range('E44').Value = '  -1.33%  '
$ws.Range('D45').Value = '''2.578'
$ws.Range('E45').Value = '  +2.27%  '
$ws.Range('D46').Value = '''4.213'
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('D48').Value = '''140.51'
$ws.Range('E48').Value = '  -1.17%  '
$ws.Range('D49').Value = '''1.323'
$ws.Range('E49').Value = '  -4.47%  '
$ws.Range('D50').Value = '''91.11'
$ws.Range('E50').Value = '  +3.14%  '
$ws.Range('D51').Value = '''0.07994'
$ws.Range('E51').Value = '  -0.64%  '
